$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "306.31"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.12"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.18%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.104"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.95%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08112"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.68%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.946"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.91%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.771"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.36%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9319"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.12%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1418"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "22.38%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1925"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "5.26%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09204"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.87%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03529"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.11%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09862"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.06%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001420"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.19%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005831"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.22%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.596"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.84%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.189"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "4.38%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.973"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.66%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.08%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1323"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.21%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.883"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.98%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2410"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.48%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.17%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004883"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.74%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001241"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.70%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02006"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.89%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04923"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.93%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01093"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "14.38%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007655"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.68%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1383"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.62%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002102"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.37%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01067"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.10%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006452"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.43%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.13%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001191"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-8.66%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.13%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.13%"
